$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new column O ("item_num") on every sheet that has the 15-col
#    Measures-style layout (Measures, ID, Dems, Dates, NewVars). This shifts
#    the existing "comment" column (old O) one place right to P.
# ---------------------------------------------------------------------------
$measureSheets = @("Measures", "ID", "Dems", "Dates", "NewVars")
foreach ($name in $measureSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(15).Insert()
    $ws.Range("O1").Value = "item_num"
}

# ---------------------------------------------------------------------------
# 2. On the Measures sheet, mark every data row (2-115) as item_num = 1.
# ---------------------------------------------------------------------------
$measures = $wb.Worksheets.Item("Measures")
$measures.Range("O2:O115").Value = 1

# ---------------------------------------------------------------------------
# 3. Update the hidden _FilterDatabase defined name to include the new col.
# ---------------------------------------------------------------------------
$wb.Names.Item("Measures!_FilterDatabase").RefersTo = "=Measures!`$A`$1:`$P`$1"

# ---------------------------------------------------------------------------
# 4. View/selection bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$idSheet = $wb.Worksheets.Item("ID")
$idSheet.Activate()
$idSheet.Range("O1:O1048576").Select()

$demsSheet = $wb.Worksheets.Item("Dems")
$demsSheet.Activate()
$demsSheet.Range("O1:O1048576").Select()

$datesSheet = $wb.Worksheets.Item("Dates")
$datesSheet.Activate()
$datesSheet.Range("O1:O1048576").Select()

$newVarsSheet = $wb.Worksheets.Item("NewVars")
$newVarsSheet.Activate()
$newVarsSheet.Range("O1:O1048576").Select()

$measures.Activate()
$measures.Range("Q112").Select()
